$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-06-15 Saturday" "2024-06-16 Sunday"

Replace-Text "162÷9=" "680÷9="
Replace-Text "855÷5=" "170÷3="
Replace-Text "365÷8=" "923÷7="
Replace-Text "965÷3=" "980÷8="
Replace-Text "421÷8=" "114÷6="
Replace-Text "754÷9=" "300÷9="
Replace-Text "175÷6=" "517÷7="
Replace-Text "374÷9=" "332÷5="
Replace-Text "572÷8=" "244÷8="
Replace-Text "495÷9=" "463÷4="
Replace-Text "567÷5=" "587÷7="
Replace-Text "196÷8=" "370÷9="
Replace-Text "744÷9=" "444÷6="
Replace-Text "810÷9=" "273÷2="
Replace-Text "341÷7=" "372÷9="
Replace-Text "928÷6=" "593÷5="
Replace-Text "838÷9=" "592÷2="
Replace-Text "778÷3=" "799÷9="
Replace-Text "682÷8=" "935÷5="
Replace-Text "489÷8=" "742÷7="
Replace-Text "664÷5=" "510÷4="
Replace-Text "895÷2=" "724÷5="
Replace-Text "651÷7=" "361÷2="
Replace-Text "169÷7=" "906÷9="
Replace-Text "116÷5=" "226÷9="
